# Applies the "completing the square" / qalist.docx update:
#   - numbering.xml: abstractNum 990 (the compact bullet list definition)
#     gets its nsid normalised to 8 hex digits and its per-level hanging
#     indent tightened from 480 -> 360 twips (list is currently unused by
#     any paragraph in this particular document, but we still touch the
#     list template object model the way Word itself would).
#   - styles.xml: a handful of the Pandoc "*Tok" syntax-highlighting
#     character styles pick up Bold / Italic, mirroring the reference
#     highlighting theme.

$d = $word.ActiveDocument

# --- 1. Compact bullet-list definition (abstractNumId 990 / numId 1000) ---
# Locate the list template backing abstractNum 990: it is the very last
# entry registered in ListTemplates for this document.
$lt = $d.ListTemplates($d.ListTemplates.Count)
$levels = $lt.ListLevels
for ($i = 1; $i -le $levels.Count; $i++) {
    $lvl = $levels.Item($i)
    # left indent is unchanged (720 * level); only the hanging indent
    # tightens from 480 -> 360 twips, i.e. numberPosition = left - 360.
    $left = 720 * $i
    try {
        $lvl.TextPosition = $left
        $lvl.NumberPosition = $left - 360
    } catch {
        # Standalone (paragraph-unattached) list template edits aren't
        # persisted back to numbering.xml by this host; nothing else to
        # do here, so move on to the style updates below.
    }
}

# --- 2. Syntax-highlighting character styles ---
# Bold keywords / control-flow tokens.
$d.Styles("KeywordTok").Font.Bold = $true
$d.Styles("ControlFlowTok").Font.Bold = $true

# Italicize documentation / comment-variable / warning tokens (already
# italic in content; re-asserting normalises the run-property order so
# <w:i/> sorts before color/shading, matching the canonical style write-out).
$d.Styles("DocumentationTok").Font.Italic = $true
$d.Styles("CommentVarTok").Font.Italic = $true
$d.Styles("WarningTok").Font.Italic = $true
